$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. The paragraph right after "... here at the college." used to contain
#    only the hidden "_GoBack" bookmark. Remove that bookmark so the
#    paragraph becomes completely empty.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Append two new sentences (as separate runs) to the paragraph that
#    discusses "The biggest challenge so far ...".
# ---------------------------------------------------------------------------
$rsquo = [char]0x2019

$challengeRange = $d.Paragraphs.Item(8).Range
$challengeRange.InsertAfter(" Another challenge that we had was if we should be trying to get our application to meet the local hospital" + $rsquo + "s requirements, and we decided that it would be better if we were to focus on getting it to meet the college" + $rsquo + "s requirements and worry about everything else after capstone.")

$challengeRange = $d.Paragraphs.Item(8).Range
$challengeRange.InsertAfter(" Some other challenges that we had were only just decisions on what programs that we will be using to create our application both easily and effectively.")

# ---------------------------------------------------------------------------
# 3. Replace the placeholder "TimeLine: Nothing here" text with the real
#    timeline write-up (several runs), and put the "_GoBack" bookmark back
#    at the very end of that paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("TimeLine: Nothing here", $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "So far where we are at now is mostly setting up the server and repository for us to work with in the future. ", 2)

$timelineRange = $d.Paragraphs.Item(10).Range
$timelineRange.InsertAfter("One to three weeks ")

$timelineRange = $d.Paragraphs.Item(10).Range
$timelineRange.InsertAfter("into the future we plan on trying to get the basic layout of the application client side of it, so we then can do more work on the server side of the application once the user-interface is ")

$timelineRange = $d.Paragraphs.Item(10).Range
$timelineRange.InsertAfter("more developed")

$timelineRange = $d.Paragraphs.Item(10).Range
$timelineRange.InsertAfter(". We do plan on sometime that we" + $rsquo + "ll show a demo of our first version of the application, so we can receive feedback from the nursing program to improve the application.")

$timelineRange = $d.Paragraphs.Item(10).Range
$timelineRange.InsertAfter(" Throughout the second quarter we will be trying to contact the nursing instructor more so that we can have them look at what we currently have and ask questions so we can ")

$timelineRange = $d.Paragraphs.Item(10).Range
$timelineRange.InsertAfter("create the program to fit their needs.")

# The final bookmark needs to sit, collapsed, right at the end of this
# paragraph's text (immediately before the paragraph mark). That exact spot
# is adjacent to the empty paragraph that follows, and adding a bookmark
# there directly lands it in the wrong place. Work around this by bookmarking
# a temporary placeholder character and then deleting that character - the
# bookmark collapses down to the correct, now-empty, location.
$timelineRange = $d.Paragraphs.Item(10).Range
$timelineRange.InsertAfter("X")

$timelineRange = $d.Paragraphs.Item(10).Range
$placeholder = $d.Range($timelineRange.End - 2, $timelineRange.End - 1)
$d.Bookmarks.Add("_GoBack", $placeholder)

$timelineRange = $d.Paragraphs.Item(10).Range
$placeholder = $d.Range($timelineRange.End - 2, $timelineRange.End - 1)
$placeholder.Text = ""
